$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (to become what was previously row 3)
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2
$f2 = $ws.Range("F2").Value2
$g2 = $ws.Range("G2").Value2
$h2 = $ws.Range("H2").Value2
$q2 = $ws.Range("Q2").Value2
$r2 = $ws.Range("R2").Value2

# Row 3 values (to become what was previously row 2)
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2
$f3 = $ws.Range("F3").Value2
$g3 = $ws.Range("G3").Value2
$h3 = $ws.Range("H3").Value2
$q3 = $ws.Range("Q3").Value2
$r3 = $ws.Range("R3").Value2

# Write swapped values into row 2
$ws.Range("A2").Value2 = $a3
$ws.Range("B2").Value2 = $b3
$ws.Range("D2").Value2 = $d3
$ws.Range("E2").Value2 = $e3
$ws.Range("F2").Value2 = $f3
$ws.Range("G2").Value2 = $g3
$ws.Range("H2").Value2 = $h3
$ws.Range("Q2").Value2 = $q3
$ws.Range("R2").Value2 = $r3

# Write swapped values into row 3
$ws.Range("A3").Value2 = $a2
$ws.Range("B3").Value2 = $b2
$ws.Range("D3").Value2 = $d2
$ws.Range("E3").Value2 = $e2
$ws.Range("F3").Value2 = $f2
$ws.Range("G3").Value2 = $g2
$ws.Range("H3").Value2 = $h2
$ws.Range("Q3").Value2 = $q2
$ws.Range("R3").Value2 = $r2
